# "Update - Task Schedular."
# The scheduled task repeatedly stamps the User sheet's A2 cell with a new
# "User_On_<dd/MM/yy>-<HH:mm>" login marker. Replay each scheduled run in
# order so the shared-string table accumulates one new unique entry per run,
# finishing with the most recent timestamp left in A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$timestamps = @(
    "User_On_19/02/19-12:30",
    "User_On_20/02/19-14:43",
    "User_On_20/02/19-15:00",
    "User_On_20/02/19-15:44",
    "User_On_25/02/19-17:14",
    "User_On_26/02/19-12:20",
    "User_On_28/02/19-15:57",
    "User_On_28/02/19-17:19",
    "User_On_28/02/19-17:47",
    "User_On_28/02/19-18:18",
    "User_On_01/03/19-11:33",
    "User_On_01/03/19-12:50",
    "User_On_01/03/19-13:06",
    "User_On_01/03/19-15:05",
    "User_On_01/03/19-15:25",
    "User_On_01/03/19-15:56",
    "User_On_08/03/19-12:22",
    "User_On_08/03/19-13:11",
    "User_On_08/03/19-14:20",
    "User_On_08/03/19-14:42"
)

foreach ($stamp in $timestamps) {
    $ws.Cells.Item(2, 1).Value = $stamp
}
